# Update project workhours spreadsheet:
# add two new logged-hours entries (10.8.2019 and 11.8.2019) to the
# bottom of the tuntikirjanpito (work-hours log) table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: 10.8.2019 entry
$ws.Range("A16").Value = "10.8.2019"
$ws.Range("B16").Value = 5
$ws.Range("C16").Value = "Käyttäjän lisäämisen ja loginin tekoa, react-bootstrap formien ja typescript tyyppien kanssa taistelua"

# Row 17: 11.8.2019 entry
$ws.Range("A17").Value = "11.8.2019"
$ws.Range("B17").Value = 5
$ws.Range("C17").Value = "Käyttäjän lisääminen ja login viimeistelty"

# Match row heights from the authored workbook
$ws.Rows.Item(16).RowHeight = 39.3
$ws.Rows.Item(17).RowHeight = 14

# Update selection/scroll position to reflect where the user ended up editing
$ws.Range("C16").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
